$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2994946666666667
$ws.Range("H2").Value = 0.8984840000000001
$ws.Range("I2").Value = 0.4989451716962827
$ws.Range("J2").Value = 0.4989451716962828
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 22.23552099999999
$ws.Range("N2").Value = 66.70656299999999
$ws.Range("O2").Value = 0.4229672982130729
$ws.Range("P2").Value = 0.422967298213073
$ws.Range("Q2").Value = 6.659419950054666
$ws.Range("R2").Value = 59.934779550492
$ws.Range("S2").Value = 0.2110374912288345
$ws.Range("T2").Value = 0.2110374912288346

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2994946666666667
$ws.Range("H3").Value = 0.8984840000000001
$ws.Range("I3").Value = 0.4989451716962827
$ws.Range("J3").Value = 0.4989451716962828
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 22.25678066666667
$ws.Range("N3").Value = 66.770342
$ws.Range("O3").Value = 0.4233717026689395
$ws.Range("P3").Value = 0.4233717026689395
$ws.Range("Q3").Value = 6.665787106836445
$ws.Range("R3").Value = 59.992083961528
$ws.Range("S3").Value = 0.2112392668795016
$ws.Range("T3").Value = 0.2112392668795016

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2994946666666667
$ws.Range("H4").Value = 0.8984840000000001
$ws.Range("I4").Value = 0.4989451716962827
$ws.Range("J4").Value = 0.4989451716962828
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.078006
$ws.Range("N4").Value = 24.234018
$ws.Range("O4").Value = 0.1536609991179876
$ws.Range("P4").Value = 0.1536609991179876
$ws.Range("Q4").Value = 2.419319714301333
$ws.Range("R4").Value = 21.773877428712
$ws.Range("S4").Value = 0.07666841358794665
$ws.Range("T4").Value = 0.07666841358794665

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.300761
$ws.Range("H5").Value = 0.902283
$ws.Range("I5").Value = 0.5010548283037172
$ws.Range("J5").Value = 0.5010548283037172
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 22.23552099999999
$ws.Range("N5").Value = 66.70656299999999
$ws.Range("O5").Value = 0.4229672982130729
$ws.Range("P5").Value = 0.422967298213073
$ws.Range("Q5").Value = 6.687577531480999
$ws.Range("R5").Value = 60.18819778332899
$ws.Range("S5").Value = 0.2119298069842384
$ws.Range("T5").Value = 0.2119298069842384

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.300761
$ws.Range("H6").Value = 0.902283
$ws.Range("I6").Value = 0.5010548283037172
$ws.Range("J6").Value = 0.5010548283037172
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 22.25678066666667
$ws.Range("N6").Value = 66.770342
$ws.Range("O6").Value = 0.4233717026689395
$ws.Range("P6").Value = 0.4233717026689395
$ws.Range("Q6").Value = 6.693971610087334
$ws.Range("R6").Value = 60.245744490786
$ws.Range("S6").Value = 0.2121324357894379
$ws.Range("T6").Value = 0.2121324357894379

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.300761
$ws.Range("H7").Value = 0.902283
$ws.Range("I7").Value = 0.5010548283037172
$ws.Range("J7").Value = 0.5010548283037172
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.078006
$ws.Range("N7").Value = 24.234018
$ws.Range("O7").Value = 0.1536609991179876
$ws.Range("P7").Value = 0.1536609991179876
$ws.Range("Q7").Value = 2.429549162566
$ws.Range("R7").Value = 21.865942463094
$ws.Range("S7").Value = 0.07699258553004089
$ws.Range("T7").Value = 0.07699258553004089

